$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$changes = @(
    @{ Cell = "H53"; Value = 690.25 }
    @{ Cell = "I53"; Value = 690.25 }
    @{ Cell = "K53"; Value = 690.25 }
    @{ Cell = "M53"; Value = -53.25 }
    @{ Cell = "H138"; Value = 312374.66 }
    @{ Cell = "J138"; Value = 442114.7 }
    @{ Cell = "L138"; Value = 1326344.1 }
    @{ Cell = "N138"; Value = -1336624.1 }
    @{ Cell = "H141"; Value = 4032.5557 }
    @{ Cell = "I141"; Value = 3477.6 }
    @{ Cell = "K141"; Value = 10432.8 }
    @{ Cell = "M141"; Value = -5252.799999999999 }
)
foreach ($c in $changes) {
    $ws.Range($c.Cell).Value = $c.Value
}

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$changes = @(
    @{ Cell = "H32"; Value = 3604.5134 }
    @{ Cell = "I32"; Value = 2853.3572 }
    @{ Cell = "K32"; Value = 2853.3572 }
    @{ Cell = "M32"; Value = -2566.3572 }
    @{ Cell = "H43"; Value = 48331.332 }
    @{ Cell = "J43"; Value = 48331.332 }
    @{ Cell = "L43"; Value = 48331.332 }
    @{ Cell = "N43"; Value = -48957.332 }
    @{ Cell = "H61"; Value = 5029.5713 }
    @{ Cell = "I61"; Value = 1999.8334 }
    @{ Cell = "J61"; Value = 7301.875 }
    @{ Cell = "K61"; Value = 1999.8334 }
    @{ Cell = "L61"; Value = 7301.875 }
    @{ Cell = "M61"; Value = -1787.8334 }
    @{ Cell = "N61"; Value = -7725.875 }
    @{ Cell = "H63"; Value = 2848.625 }
    @{ Cell = "I63"; Value = 2848.625 }
    @{ Cell = "K63"; Value = 2848.625 }
    @{ Cell = "M63"; Value = -2162.625 }
    @{ Cell = "H66"; Value = 2848.625 }
    @{ Cell = "I66"; Value = 2848.625 }
    @{ Cell = "K66"; Value = 14243.125 }
    @{ Cell = "M66"; Value = -10811.125 }
    @{ Cell = "H136"; Value = 5029.5713 }
    @{ Cell = "I136"; Value = 1999.8334 }
    @{ Cell = "J136"; Value = 7301.875 }
    @{ Cell = "K136"; Value = 5999.5002 }
    @{ Cell = "L136"; Value = 21905.625 }
    @{ Cell = "M136"; Value = -3449.5002 }
    @{ Cell = "N136"; Value = -27005.625 }
)
foreach ($c in $changes) {
    $ws.Range($c.Cell).Value = $c.Value
}

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$changes = @(
    @{ Cell = "H107"; Value = 4050218 }
    @{ Cell = "I107"; Value = 4526389 }
    @{ Cell = "J107"; Value = 2766.5 }
    @{ Cell = "K107"; Value = 4526389 }
    @{ Cell = "L107"; Value = 2766.5 }
    @{ Cell = "M107"; Value = -4524469 }
    @{ Cell = "N107"; Value = -6606.5 }
    @{ Cell = "H134"; Value = 2058.4092 }
    @{ Cell = "I134"; Value = 1575.5883 }
    @{ Cell = "K134"; Value = 4726.7649 }
    @{ Cell = "M134"; Value = -2191.7649 }
)
foreach ($c in $changes) {
    $ws.Range($c.Cell).Value = $c.Value
}

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$changes = @(
    @{ Cell = "H16"; Value = 2338.611 }
    @{ Cell = "I16"; Value = 2311.3845 }
    @{ Cell = "K16"; Value = 2311.3845 }
    @{ Cell = "M16"; Value = -2024.3845 }
    @{ Cell = "H31"; Value = 4771.1934 }
    @{ Cell = "I31"; Value = 4588 }
    @{ Cell = "K31"; Value = 4588 }
    @{ Cell = "M31"; Value = -4293 }
    @{ Cell = "H34"; Value = 4771.1934 }
    @{ Cell = "I34"; Value = 4588 }
    @{ Cell = "K34"; Value = 4588 }
    @{ Cell = "M34"; Value = -4386 }
    @{ Cell = "H58"; Value = 3545.9565 }
    @{ Cell = "I58"; Value = 2608.25 }
    @{ Cell = "K58"; Value = 2608.25 }
    @{ Cell = "M58"; Value = -2405.25 }
    @{ Cell = "H113"; Value = 2338.611 }
    @{ Cell = "I113"; Value = 2311.3845 }
    @{ Cell = "K113"; Value = 2311.3845 }
    @{ Cell = "M113"; Value = -141.3845000000001 }
    @{ Cell = "H132"; Value = 26320064 }
    @{ Cell = "I132"; Value = 38464940 }
    @{ Cell = "K132"; Value = 115394820 }
    @{ Cell = "M132"; Value = -115392290 }
    @{ Cell = "H136"; Value = 3545.9565 }
    @{ Cell = "I136"; Value = 2608.25 }
    @{ Cell = "K136"; Value = 7824.75 }
    @{ Cell = "M136"; Value = -5274.75 }
)
foreach ($c in $changes) {
    $ws.Range($c.Cell).Value = $c.Value
}

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$changes = @(
    @{ Cell = "H68"; Value = 14290729 }
    @{ Cell = "J68"; Value = 16668351 }
    @{ Cell = "L68"; Value = 50005053 }
    @{ Cell = "N68"; Value = -50006675 }
    @{ Cell = "H71"; Value = 14290729 }
    @{ Cell = "J71"; Value = 16668351 }
    @{ Cell = "L71"; Value = 150015159 }
    @{ Cell = "N71"; Value = -150023271 }
    @{ Cell = "H80"; Value = 4358 }
    @{ Cell = "I80"; Value = 5072.6665 }
    @{ Cell = "J80"; Value = 3929.2 }
    @{ Cell = "K80"; Value = 15217.9995 }
    @{ Cell = "L80"; Value = 11787.6 }
    @{ Cell = "M80"; Value = -14281.9995 }
    @{ Cell = "N80"; Value = -13659.6 }
    @{ Cell = "H83"; Value = 4358 }
    @{ Cell = "I83"; Value = 5072.6665 }
    @{ Cell = "J83"; Value = 3929.2 }
    @{ Cell = "K83"; Value = 45653.9985 }
    @{ Cell = "L83"; Value = 35362.8 }
    @{ Cell = "M83"; Value = -40973.9985 }
    @{ Cell = "N83"; Value = -44722.8 }
    @{ Cell = "H92"; Value = 464.7143 }
    @{ Cell = "J92"; Value = 375.66666 }
    @{ Cell = "L92"; Value = 1126.99998 }
    @{ Cell = "N92"; Value = -3622.99998 }
    @{ Cell = "H109"; Value = 1875.3636 }
    @{ Cell = "J109"; Value = 3763.75 }
    @{ Cell = "L109"; Value = 11291.25 }
    @{ Cell = "N109"; Value = -13371.25 }
    @{ Cell = "H113"; Value = 6044.273 }
    @{ Cell = "I113"; Value = 496 }
    @{ Cell = "K113"; Value = 1488 }
    @{ Cell = "M113"; Value = 682 }
    @{ Cell = "H122"; Value = 892.46155 }
    @{ Cell = "I122"; Value = 748.5714 }
    @{ Cell = "J122"; Value = 945.4737 }
    @{ Cell = "K122"; Value = 6737.1426 }
    @{ Cell = "L122"; Value = 8509.263300000001 }
    @{ Cell = "M122"; Value = -4287.1426 }
    @{ Cell = "N122"; Value = -13409.2633 }
    @{ Cell = "H129"; Value = 1463.4667 }
    @{ Cell = "I129"; Value = 894.1667 }
    @{ Cell = "K129"; Value = 2682.5001 }
    @{ Cell = "M129"; Value = 2317.4999 }
    @{ Cell = "H137"; Value = 3423.9333 }
    @{ Cell = "I137"; Value = 3517.3845 }
    @{ Cell = "K137"; Value = 10552.1535 }
    @{ Cell = "M137"; Value = -5452.1535 }
    @{ Cell = "H138"; Value = 3174.2666 }
    @{ Cell = "I138"; Value = 1859.5834 }
    @{ Cell = "K138"; Value = 5578.7502 }
    @{ Cell = "M138"; Value = -438.7502000000004 }
)
foreach ($c in $changes) {
    $ws.Range($c.Cell).Value = $c.Value
}

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$changes = @(
    @{ Cell = "H113"; Value = 4098.8066 }
    @{ Cell = "I113"; Value = 4110.4 }
    @{ Cell = "J113"; Value = 4050.5 }
    @{ Cell = "K113"; Value = 4110.4 }
    @{ Cell = "L113"; Value = 4050.5 }
    @{ Cell = "M113"; Value = -1940.4 }
    @{ Cell = "N113"; Value = -8390.5 }
)
foreach ($c in $changes) {
    $ws.Range($c.Cell).Value = $c.Value
}

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$changes = @(
    @{ Cell = "H136"; Value = 5610.647 }
    @{ Cell = "I136"; Value = 5762.909 }
    @{ Cell = "J136"; Value = 5331.5 }
    @{ Cell = "K136"; Value = 17288.727 }
    @{ Cell = "L136"; Value = 15994.5 }
    @{ Cell = "M136"; Value = -14738.727 }
    @{ Cell = "N136"; Value = -21094.5 }
)
foreach ($c in $changes) {
    $ws.Range($c.Cell).Value = $c.Value
}

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$changes = @(
    @{ Cell = "H132"; Value = 9262387 }
    @{ Cell = "I132"; Value = 11114165 }
    @{ Cell = "J132"; Value = 3499.8333 }
    @{ Cell = "K132"; Value = 33342495 }
    @{ Cell = "L132"; Value = 10499.4999 }
    @{ Cell = "M132"; Value = -33339965 }
    @{ Cell = "N132"; Value = -15559.4999 }
    @{ Cell = "H136"; Value = 28574232 }
    @{ Cell = "I136"; Value = 33334576 }
    @{ Cell = "J136"; Value = 12168 }
    @{ Cell = "K136"; Value = 100003728 }
    @{ Cell = "L136"; Value = 36504 }
    @{ Cell = "M136"; Value = -100001178 }
    @{ Cell = "N136"; Value = -41604 }
)
foreach ($c in $changes) {
    $ws.Range($c.Cell).Value = $c.Value
}
